# Update "想去人数" (F column) values on sheets "展览" and "全部类型"
# Rows 2-5: F2=46, F3=274, F4=17, F5=61

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
$values = @{ 2 = 46; 3 = 274; 4 = 17; 5 = 61 }

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($row in $values.Keys) {
        $ws.Cells.Item($row, 6).Value = $values[$row]
    }
}
